$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New identifiers / timestamps produced by this handback run
# ---------------------------------------------------------------------------
$guidZh   = "6808ca51-4609-41ca-bd42-6a75a51e4a3e"
$guidDe   = "ffff00a57b2b-3c4a-4f17-b157-9ebacecac0b7"
$xlfHash  = "ae41ee51ac863b8283b38ce29343d0bb75a8253e"

$mdZh     = "$guidZh.md"
$mdDe     = "$guidDe.md"
$pathZh   = "e2e\$mdZh"
$pathDe   = "e2e\$mdDe"

$xlfZhCn  = "$guidZh.$xlfHash.zh-cn.xlf"
$xlfDeDe  = "$guidZh.$xlfHash.de-de.xlf"

$latestHoDate = "2016-09-06 11:23:10"
$zhGenDate    = "2016-09-06 11:22:57"
$zhBackDate   = "2016-09-06 11:23:32"
$deGenDate    = "$latestHoDate"
$deBackDate   = "2016-09-06 11:23:41"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = $mdZh
$ws1.Range("B2").Value = $pathZh
$ws1.Range("G2").Value = $latestHoDate

$ws1.Range("A3").Value = $mdDe
$ws1.Range("B3").Value = $pathDe
$ws1.Range("G3").Value = $latestHoDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/0e0851df-7c52-45e6-b4b4-094126f96cec.md", "", "", $pathZh)
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/6f9b1139-1bbe-4751-a6b3-7ede80dff354.md", "", "", $pathDe)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = $mdZh
$ws2.Range("G2").Value = $xlfZhCn
$ws2.Range("H2").Value = $zhGenDate
$ws2.Range("I2").Value = $mdZh
$ws2.Range("J2").Value = $xlfZhCn
$ws2.Range("K2").Value = $zhBackDate

$ws2.Range("A3").Value = $mdDe
$ws2.Range("G3").Value = $xlfZhCn
$ws2.Range("H3").Value = $zhGenDate
$ws2.Range("I3").Value = $mdDe
$ws2.Range("J3").Value = $xlfZhCn
$ws2.Range("K3").Value = $zhBackDate

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/0e0851df-7c52-45e6-b4b4-094126f96cec.md", "", "", $mdZh)
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ad94a87ec331acfaba7f3a2281209ec5c937178b/e2e/0e0851df-7c52-45e6-b4b4-094126f96cec.md", "", "", $mdZh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/6f9b1139-1bbe-4751-a6b3-7ede80dff354.md", "", "", $mdDe)
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/ad94a87ec331acfaba7f3a2281209ec5c937178b/e2e/6f9b1139-1bbe-4751-a6b3-7ede80dff354.md", "", "", $mdDe)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = $mdZh
$ws3.Range("G2").Value = $xlfDeDe
$ws3.Range("H2").Value = $deGenDate
$ws3.Range("I2").Value = $mdZh
$ws3.Range("J2").Value = $xlfDeDe
$ws3.Range("K2").Value = $deBackDate

$ws3.Range("A3").Value = $mdDe
$ws3.Range("G3").Value = $xlfDeDe
$ws3.Range("H3").Value = $deGenDate
$ws3.Range("I3").Value = $mdDe
$ws3.Range("J3").Value = $xlfDeDe
$ws3.Range("K3").Value = $deBackDate

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/0e0851df-7c52-45e6-b4b4-094126f96cec.md", "", "", $mdZh)
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/29349bd427f478a45b1866862bf985f027f3836a/e2e/0e0851df-7c52-45e6-b4b4-094126f96cec.md", "", "", $mdZh)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6859e6130c153109bc151a2207a4f67a9d445c64/e2e/6f9b1139-1bbe-4751-a6b3-7ede80dff354.md", "", "", $mdDe)
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/29349bd427f478a45b1866862bf985f027f3836a/e2e/6f9b1139-1bbe-4751-a6b3-7ede80dff354.md", "", "", $mdDe)
